$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rubric row 11 "Dimensions of diversity in the community" is being re-weighted
#     from 10 points down to 5 points, and gets its own new point scale (5-4 / 3-2 / 1-0)
#     instead of sharing the common 10-8 / 7-4 / 3-0 scale ---
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = "5–4"
$ws.Range("F11").Value = "3–2"
$ws.Range("G11").Value = "1–0"

# Nudge the number format on E11 so it picks up its own (new) cell style, matching the
# freshly-created point scale rather than continuing to share style with the rows above/below.
$ws.Range("E11").NumberFormat = "General"

# --- Rubric row 17 "Implementation" is being re-weighted from 15 points down to 10 points,
#     so it now shares the common 10-8 / 7-4 / 3-0 point scale used elsewhere ---
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = "10–8"
$ws.Range("F17").Value = "7–4"
$ws.Range("G17").Value = "3–0"

# Recalculate the Subtotal/Total formulas that depend on D11 and D17
$wb.Application.Calculate()
